$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# The first three rows' values get replaced with the placeholder "0M",
# and their original values move down to become the sole content of the
# last three rows (which previously held tab-separated multi-value runs).
$origRow1 = $t.Rows.Item(1).Cells.Item(1).Range.Text
$origRow2 = $t.Rows.Item(2).Cells.Item(1).Range.Text
$origRow3 = $t.Rows.Item(3).Cells.Item(1).Range.Text

$t.Rows.Item(1).Cells.Item(1).Range.Text = "0M"
$t.Rows.Item(2).Cells.Item(1).Range.Text = "0M"
$t.Rows.Item(3).Cells.Item(1).Range.Text = "0M"

# Insert 10 new single-value rows right after row 3 (before the old row 4).
$newValues = @("439", "0.00003", "0.00014", "0.00004", "0.00001", "0.00004", "0.00005", "0.00005", "0.02022", "100.0")
$insertPos = 4
foreach ($v in $newValues) {
    $refRow = $t.Rows.Item($insertPos)
    $newRow = $t.Rows.Add($refRow)
    $newRow.Cells.Item(1).Range.Text = $v
    $insertPos = $insertPos + 1
}

# Collapse the (now shifted) final three rows, which held tab-separated
# multi-run values, down to a single run holding the values that used to
# be in rows 1-3.
$rowCount = $t.Rows.Count
$t.Rows.Item($rowCount - 2).Cells.Item(1).Range.Text = $origRow1
$t.Rows.Item($rowCount - 1).Cells.Item(1).Range.Text = $origRow2
$t.Rows.Item($rowCount).Cells.Item(1).Range.Text = $origRow3
